$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (partial run edits to preserve rich-text formatting) ---
# A8: "Volume 30   Number  31" -> "Volume 30   Number  32"
$a8 = $ws.Range("A8")
$a8.Characters(21,2).Text = "32"

# C9: "Report Covering the Week  7/31/2023  Through  8/6/2023"
#  -> "Report Covering the Week  8/7/2023  Through  8/13/2023"
# Replace the right-most substring first so the left substring's character
# offsets aren't shifted by the differing replacement length.
$c9 = $ws.Range("C9")
$c9.Characters(47,8).Text = "8/13/2023"
$c9.Characters(27,9).Text = "8/7/2023"

# --- Weekly crime-statistics table updates (rows 14-30) ---
$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 23
$ws.Range("G14").Value = 35
$ws.Range("H14").Value = -34.285714285714
$ws.Range("I14").Value = 247
$ws.Range("J14").Value = 280
$ws.Range("K14").Value = -11.785714285714
$ws.Range("L14").Value = -15.699658703071
$ws.Range("M14").Value = -27.138643067846
$ws.Range("N14").Value = -79.295892707460
$ws.Range("C15").Value = 25
$ws.Range("D15").Value = 29
$ws.Range("E15").Value = -13.793103448275
$ws.Range("F15").Value = 103
$ws.Range("G15").Value = 139
$ws.Range("H15").Value = -25.899280575539
$ws.Range("I15").Value = 897
$ws.Range("J15").Value = 1025
$ws.Range("K15").Value = -12.487804878048
$ws.Range("L15").Value = -3.961456102783
$ws.Range("M15").Value = 11.290322580645
$ws.Range("N15").Value = -55.528011898859
$ws.Range("C16").Value = 335
$ws.Range("D16").Value = 391
$ws.Range("E16").Value = -14.322250639386
$ws.Range("F16").Value = 1417
$ws.Range("G16").Value = 1560
$ws.Range("H16").Value = -9.166666666666
$ws.Range("I16").Value = 9976
$ws.Range("J16").Value = 10652
$ws.Range("K16").Value = -6.346226060833
$ws.Range("L16").Value = 31.280431635741
$ws.Range("M16").Value = -11.662091561144
$ws.Range("N16").Value = -80.644910946412
$ws.Range("C17").Value = 516
$ws.Range("D17").Value = 556
$ws.Range("E17").Value = -7.194244604316
$ws.Range("F17").Value = 2318
$ws.Range("G17").Value = 2316
$ws.Range("H17").Value = 0.086355785837
$ws.Range("I17").Value = 17150
$ws.Range("J17").Value = 16290
$ws.Range("K17").Value = 5.279312461632
$ws.Range("L17").Value = 26.896041435442
$ws.Range("M17").Value = 60.821455363841
$ws.Range("N17").Value = -33.753090234857
$ws.Range("C18").Value = 237
$ws.Range("D18").Value = 305
$ws.Range("E18").Value = -22.295081967213
$ws.Range("F18").Value = 1088
$ws.Range("G18").Value = 1215
$ws.Range("H18").Value = -10.452674897119
$ws.Range("I18").Value = 8554
$ws.Range("J18").Value = 9551
$ws.Range("K18").Value = -10.438697518584
$ws.Range("L18").Value = 19.003895381190
$ws.Range("M18").Value = -22.504076825511
$ws.Range("N18").Value = -86.018078098694
$ws.Range("C19").Value = 1034
$ws.Range("D19").Value = 1077
$ws.Range("E19").Value = -3.992571959145
$ws.Range("F19").Value = 4119
$ws.Range("G19").Value = 4338
$ws.Range("H19").Value = -5.048409405255
$ws.Range("I19").Value = 30567
$ws.Range("J19").Value = 31312
$ws.Range("K19").Value = -2.379279509453
$ws.Range("L19").Value = 44.812393405343
$ws.Range("M19").Value = 36.234790747426
$ws.Range("N19").Value = -40.990347490347
$ws.Range("C20").Value = 357
$ws.Range("D20").Value = 258
$ws.Range("E20").Value = 38.372093023255
$ws.Range("F20").Value = 1415
$ws.Range("G20").Value = 1067
$ws.Range("H20").Value = 32.614807872539
$ws.Range("I20").Value = 9658
$ws.Range("J20").Value = 8144
$ws.Range("K20").Value = 18.590373280943
$ws.Range("L20").Value = 67.383015597920
$ws.Range("M20").Value = 52.719797596457
$ws.Range("N20").Value = -85.845973474023
$ws.Range("C21").Value = 2509
$ws.Range("D21").Value = 2621
$ws.Range("E21").Value = -4.273178176268
$ws.Range("F21").Value = 10483
$ws.Range("G21").Value = 10670
$ws.Range("H21").Value = -1.752577319587
$ws.Range("I21").Value = 77049
$ws.Range("J21").Value = 77254
$ws.Range("K21").Value = -0.265358428042
$ws.Range("L21").Value = 36.594748878685
$ws.Range("M21").Value = 22.492488195736
$ws.Range("N21").Value = -70.575587923041
$ws.Range("C22").Value = 35
$ws.Range("D22").Value = 42
$ws.Range("E22").Value = -16.666666666666
$ws.Range("F22").Value = 148
$ws.Range("G22").Value = 147
$ws.Range("H22").Value = 0.680272108843
$ws.Range("I22").Value = 1340
$ws.Range("J22").Value = 1400
$ws.Range("K22").Value = -4.285714285714
$ws.Range("L22").Value = 45.336225596529
$ws.Range("M22").Value = 3.795507358636
$ws.Range("C23").Value = 137
$ws.Range("D23").Value = 141
$ws.Range("E23").Value = -2.836879432624
$ws.Range("F23").Value = 520
$ws.Range("G23").Value = 537
$ws.Range("H23").Value = -3.165735567970
$ws.Range("I23").Value = 3855
$ws.Range("J23").Value = 3737
$ws.Range("K23").Value = 3.157613058603
$ws.Range("L23").Value = 17.602196461256
$ws.Range("M23").Value = 51.295133437990
$ws.Range("C24").Value = 2295
$ws.Range("D24").Value = 2323
$ws.Range("E24").Value = -1.205337925096
$ws.Range("F24").Value = 9065
$ws.Range("G24").Value = 9521
$ws.Range("H24").Value = -4.789412876798
$ws.Range("I24").Value = 67811
$ws.Range("J24").Value = 69785
$ws.Range("K24").Value = -2.828688113491
$ws.Range("L24").Value = 38.468921016090
$ws.Range("M24").Value = 37.592322051781
$ws.Range("C25").Value = 903
$ws.Range("D25").Value = 817
$ws.Range("E25").Value = 10.526315789473
$ws.Range("F25").Value = 3559
$ws.Range("G25").Value = 3298
$ws.Range("H25").Value = 7.913887204366
$ws.Range("I25").Value = 27048
$ws.Range("J25").Value = 25816
$ws.Range("K25").Value = 4.772234273318
$ws.Range("L25").Value = 28.689694547530
$ws.Range("M25").Value = -6.168042739193
$ws.Range("D26").Value = 47
$ws.Range("E26").Value = -10.638297872340
$ws.Range("F26").Value = 172
$ws.Range("G26").Value = 197
$ws.Range("H26").Value = -12.690355329949
$ws.Range("I26").Value = 1497
$ws.Range("J26").Value = 1646
$ws.Range("K26").Value = -9.052247873633
$ws.Range("L26").Value = -0.729442970822
$ws.Range("C27").Value = 107
$ws.Range("D27").Value = 109
$ws.Range("E27").Value = -1.834862385321
$ws.Range("F27").Value = 428
$ws.Range("G27").Value = 431
$ws.Range("H27").Value = -0.696055684454
$ws.Range("I27").Value = 3261
$ws.Range("J27").Value = 3175
$ws.Range("K27").Value = 2.708661417322
$ws.Range("L27").Value = 12.759336099585
$ws.Range("C28").Value = 16
$ws.Range("D28").Value = 36
$ws.Range("E28").Value = -55.555555555555
$ws.Range("F28").Value = 106
$ws.Range("G28").Value = 159
$ws.Range("H28").Value = -33.333333333333
$ws.Range("I28").Value = 750
$ws.Range("J28").Value = 1046
$ws.Range("K28").Value = -28.298279158699
$ws.Range("L28").Value = -34.497816593886
$ws.Range("M28").Value = -32.553956834532
$ws.Range("N28").Value = -79.685807150595
$ws.Range("C29").Value = 15
$ws.Range("D29").Value = 31
$ws.Range("E29").Value = -51.612903225806
$ws.Range("F29").Value = 90
$ws.Range("G29").Value = 129
$ws.Range("H29").Value = -30.232558139534
$ws.Range("I29").Value = 633
$ws.Range("J29").Value = 865
$ws.Range("K29").Value = -26.820809248554
$ws.Range("L29").Value = -34.876543209876
$ws.Range("M29").Value = -31.045751633986
$ws.Range("N29").Value = -80.956678700361
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 10
$ws.Range("E30").Value = -80
$ws.Range("F30").Value = 23
$ws.Range("G30").Value = 59
$ws.Range("H30").Value = -61.016949152542
$ws.Range("I30").Value = 295
$ws.Range("J30").Value = 426
$ws.Range("K30").Value = -30.751173708920
$ws.Range("L30").Value = -14.492753623188
